$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column-D price cells to retain their exact text representation
# (Excel would otherwise auto-convert numeric-looking strings to floats,
# dropping trailing zeros / precision / or flipping to scientific notation).
$priceCells = @('D2', 'D3', 'D5', 'D6', 'D8', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D33', 'D34', 'D36', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D48', 'D49', 'D50', 'D51')
foreach ($ref in $priceCells) { $ws.Range($ref).NumberFormat = "@" }

$ws.Range('D2').Value = '29.975.81'
$ws.Range('E2').Value = '  +2.82%  '
$ws.Range('D3').Value = '1.861.97'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '245.97'
$ws.Range('E5').Value = '  +1.81%  '
$ws.Range('D6').Value = '0.6389'
$ws.Range('E6').Value = '  +3.49%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '0.3000'
$ws.Range('E8').Value = '  +3.73%  '
$ws.Range('E9').Value = '  +1.88%  '
$ws.Range('D10').Value = '24.48'
$ws.Range('E10').Value = '  +6.48%  '
$ws.Range('D11').Value = '0.07684'
$ws.Range('E11').Value = '  +0.19%  '
$ws.Range('D12').Value = '1.873.29'
$ws.Range('E12').Value = '  +2.76%  '
$ws.Range('D13').Value = '5.061'
$ws.Range('E13').Value = '  +2.25%  '
$ws.Range('D14').Value = '0.6901'
$ws.Range('E14').Value = '  +4.36%  '
$ws.Range('D15').Value = '84.18'
$ws.Range('E15').Value = '  +2.77%  '
$ws.Range('D16').Value = '0.000009459'
$ws.Range('E16').Value = '  +6.14%  '
$ws.Range('D17').Value = '6.087'
$ws.Range('E17').Value = '  +4.28%  '
$ws.Range('D18').Value = '29.952.29'
$ws.Range('E18').Value = '  +2.90%  '
$ws.Range('D19').Value = '2.123.63'
$ws.Range('E19').Value = '  +2.53%  '
$ws.Range('D20').Value = '240.21'
$ws.Range('E20').Value = '  +1.45%  '
$ws.Range('D21').Value = '12.69'
$ws.Range('D22').Value = '1.0000'
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').Value = '7.426'
$ws.Range('E23').Value = '  +4.06%  '
$ws.Range('D24').Value = '1.000'
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').Value = '159.65'
$ws.Range('E25').Value = '  +1.15%  '
$ws.Range('D26').Value = '0.1424'
$ws.Range('E26').Value = '  +0.81%  '
$ws.Range('D27').Value = '8.583'
$ws.Range('E27').Value = '  +1.82%  '
$ws.Range('D28').Value = '18.01'
$ws.Range('E28').Value = '  +2.33%  '
$ws.Range('D29').Value = '0.06138'
$ws.Range('E29').Value = '  +10.63%  '
$ws.Range('D30').Value = '1.507'
$ws.Range('E30').Value = '  +1.55%  '
$ws.Range('D31').Value = '1.281'
$ws.Range('E31').Value = '  +6.56%  '
$ws.Range('E32').Value = '  +1.19%  '
$ws.Range('D33').Value = '4.148'
$ws.Range('E33').Value = '  +1.30%  '
$ws.Range('D34').Value = '1.889'
$ws.Range('E34').Value = '  +3.45%  '
$ws.Range('E35').Value = '  +3.10%  '
$ws.Range('D36').Value = '0.7360'
$ws.Range('E36').Value = '  +0.19%  '
$ws.Range('E37').Value = '  -0.12%  '
$ws.Range('D38').Value = '2.867'
$ws.Range('E38').Value = '  +1.10%  '
$ws.Range('D39').Value = '0.01804'
$ws.Range('E39').Value = '  +2.66%  '
$ws.Range('D40').Value = '1.224.25'
$ws.Range('E40').Value = '  +1.23%  '
$ws.Range('D41').Value = '0.9299'
$ws.Range('E41').Value = '  +3.28%  '
$ws.Range('D42').Value = '6.274'
$ws.Range('E42').Value = '  -0.89%  '
$ws.Range('B43').Value = 'RocketPoolETH'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D43').Value = '2.036.32'
$ws.Range('E43').Value = '  +3.15%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('D45').Value = '102.22'
$ws.Range('E45').Value = '  +0.79%  '
$ws.Range('D46').Value = '66.40'
$ws.Range('E46').Value = '  +2.76%  '
$ws.Range('E47').Value = '  +0.46%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '9.344'
$ws.Range('E48').Value = '  +3.66%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.00000000119'
$ws.Range('E49').Value = '  -7.12%  '
$ws.Range('D50').Value = '0.4101'
$ws.Range('E50').Value = '  +2.39%  '
$ws.Range('D51').Value = '0.1143'
$ws.Range('E51').Value = '  +3.35%  '
